# Update the "Fitness" values (column C) on Sheet1 to reflect the new run log values.
# Column A = Run, Column B = Generation (0-based, row-2), Column C = Fitness.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 11408
$ws.Range("C3:C5").Value = 11307
$ws.Range("C6").Value = 10537
$ws.Range("C7:C10").Value = 9202
$ws.Range("C11").Value = 9144
$ws.Range("C12:C13").Value = 9072
$ws.Range("C14:C21").Value = 8706
$ws.Range("C22:C23").Value = 8404
$ws.Range("C24").Value = 8365
$ws.Range("C25:C26").Value = 8340
$ws.Range("C27:C31").Value = 7917
$ws.Range("C32:C34").Value = 7884
$ws.Range("C35").Value = 7345
$ws.Range("C36:C50").Value = 7312
$ws.Range("C167:C252").Value = 7310
